$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new test row: Purchase Foreign Currency testcase
$ws.Range("A3").Value = "new"
$ws.Range("B3").Value = "user"

# Move selection to B3, matching the post-edit saved selection state
$ws.Range("B3").Select()
